$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALX")

# Insert two new columns before column D (old D:K data shifts to F:M)
$ws.Range("D:E").EntireColumn.Insert()

# Helper to fill a "dates" row header (style 2, numeric date values) for new D/E
function Set-DateRow($row) {
    $ws.Cells.Item($row, 4).Value = 43465
    $ws.Cells.Item($row, 5).Value = 43373
}

# Helper to set plain numeric values for new D/E on a given row
function Set-Values($row, $dVal, $eVal) {
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal
}

# ---- Income Statement (rows 7-35) ----
Set-DateRow 7
Set-Values 8 57600 59100
Set-Values 9 24200 27000
Set-Values 10 33400 32100
# row 11 stays blank
Set-Values 12 "NA" "NA"
Set-Values 13 0 0
Set-Values 14 0 0
Set-Values 15 7900 8200
# row 16 stays blank
Set-Values 17 32700 35800
Set-Values 18 24900 23300
# row 19 stays blank
Set-Values 20 -2500 3000
Set-Values 21 31600 35900
Set-Values 22 12400 11300
Set-Values 23 10000 15000
Set-Values 24 0 0
Set-Values 25 0 0
Set-Values 26 10000 15000
Set-Values 27 10000 15000
Set-Values 28 0 0
Set-Values 29 0 0
Set-Values 30 0 0
Set-Values 31 0 0
Set-Values 32 2500 -3000
Set-Values 33 10000 15000
Set-Values 34 0 0
Set-Values 35 10000 15000

# ---- Balance Sheet (rows 38-77) ----
Set-DateRow 38
# row 39 stays blank
# row 40 stays blank
Set-Values 41 283100 303700
Set-Values 42 0 0
Set-Values 43 172900 172800
Set-Values 44 0 0
Set-Values 45 0 0
Set-Values 46 0 0
Set-Values 47 218900 225800
Set-Values 48 730300 735000
Set-Values 49 0 0
Set-Values 50 0 0
Set-Values 51 0 0
Set-Values 52 47100 47300
Set-Values 53 0 0
Set-Values 54 1481300 1500900
# row 55 stays blank
# row 56 stays blank
Set-Values 57 30900 36300
Set-Values 58 0 0
Set-Values 59 700 400
Set-Values 60 0 0
Set-Values 61 1161500 1162900
Set-Values 62 "NA" "NA"
Set-Values 63 0 0
Set-Values 64 0 0
Set-Values 65 0 0
Set-Values 66 1196200 1202700
# row 67 stays blank
Set-Values 68 0 0
Set-Values 69 0 0
Set-Values 70 0 0
Set-Values 71 0 0
Set-Values 72 248400 261500
Set-Values 73 0 0
Set-Values 74 0 0
Set-Values 75 0 0
Set-Values 76 285100 298100
Set-Values 77 0 0

# ---- Cash Flow Statement (rows 80-102) ----
Set-DateRow 80
Set-Values 81 10000 15000
# row 82 stays blank
Set-Values 83 9200 9600
Set-Values 84 0 0
Set-Values 85 0 0
Set-Values 86 0 0
Set-Values 87 0 0
Set-Values 88 0 0
Set-Values 89 6900 33200
# row 90 stays blank
Set-Values 91 -1500 -700
Set-Values 92 0 0
Set-Values 93 0 0
Set-Values 94 -900 100
# row 95 stays blank
Set-Values 96 -23000 -23000
Set-Values 97 0 0
Set-Values 98 0 0
Set-Values 99 0 0
Set-Values 100 -25700 -24100
Set-Values 101 0 0
Set-Values 102 -19800 9200
